# Amgen DCF - "new data on Amgen"
# Update the Market Return assumption on the DCF Model sheet; every other
# changed cell in the workbook is a formula that recalculates from this one
# input, so a single write (plus re-positioning the selection, matching the
# author's last on-screen selection) reproduces the whole diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DCF Model")
$ws.Activate()

# B15 = "Market Return" assumption: 1.68% -> 7.00%
$ws.Range("G15").Value = 0.07

# Leave the sheet scrolled/selected where the author left it (G16).
$ws.Range("G16").Select()
